$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: copy A22's date/time number format onto A23, then set the
# new timestamp value (keeps style index "1" instead of minting a new xf).
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = 42604.890173611115

$ws.Range("B23").Value = "Bag"

$ws.Range("C23").Value = 6525
$ws.Range("D23").Value = 9801
$ws.Range("E23").Value = 1190
$ws.Range("F23").Value = 108
$ws.Range("G23").Value = 105
$ws.Range("H23").Value = 49
$ws.Range("I23").Value = 48
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = 18
$ws.Range("L23").Value = 21
$ws.Range("M23").Value = 78
